$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a date-like (or otherwise number-ish) string into a cell
# as a genuine shared-string (plain text, default style) instead of
# letting Excel auto-convert strings such as "2/13/2010" into a numeric
# date serial. We do this by building the text via a formula in a
# scratch cell (whose computed result is always a string), copying it,
# then pasting "Values only" into the destination - PasteSpecial(values)
# keeps the text as text without re-inferring a date/number type.
function Set-TextValue($cell, $text) {
    $scratch = $ws.Cells.Item(1000, 1000)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.Clear()
}

# Fill in the new status-report rows. The shared-string table is built up
# in the same order the entries were originally authored (not strictly
# top-to-bottom by row), which is why "2/14/2010" ends up duplicated in
# the string table while "2/18/2010" ends up shared between rows 42/43.
Set-TextValue $ws.Cells.Item(39,1) "2/14/2010"
Set-TextValue $ws.Cells.Item(39,3) "Encoder speed requirement calculation"

Set-TextValue $ws.Cells.Item(41,1) "2/15/2010"
Set-TextValue $ws.Cells.Item(41,3) "Ordered control system components"

Set-TextValue $ws.Cells.Item(40,3) "Researched QNX installation methods and issues"
Set-TextValue $ws.Cells.Item(40,1) "2/14/2010"

Set-TextValue $ws.Cells.Item(38,1) "2/13/2010"
Set-TextValue $ws.Cells.Item(38,3) "Researched sterilizable motor"

Set-TextValue $ws.Cells.Item(42,1) "2/18/2010"
Set-TextValue $ws.Cells.Item(42,3) "Group Meeting on Skype"

Set-TextValue $ws.Cells.Item(43,3) "Gravity Compensation motor experiment"
Set-TextValue $ws.Cells.Item(43,1) "2/18/2010"

$ws.Cells.Item(38,2).Value = 1.5
$ws.Cells.Item(39,2).Value = 0.5
$ws.Cells.Item(40,2).Value = 2
$ws.Cells.Item(41,2).Value = 3
$ws.Cells.Item(42,2).Value = 1
$ws.Cells.Item(43,2).Value = 0.5

$ws.Range("A44").Select()
